$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed coin snapshot (Coin, Link, Price, Volume(1h)) per row.
# A leading apostrophe forces plain-decimal Price strings (e.g. "1.00", "0.0000219")
# to stay text instead of being normalized to a number by Excel - matching the
# original inline-string cell content exactly.
$rows = @(
    @{ Row = 2; B = "Bitcoin"; C = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D = "67.401.97"; E = "  -1.42%  " },
    @{ Row = 3; B = "Ethereum"; C = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D = "3.505.95"; E = "  -2.51%  " },
    @{ Row = 4; B = "TetherUSD"; C = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D = "'0.999"; E = "  -0.07%  " },
    @{ Row = 5; B = "BNB"; C = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D = "'607.57"; E = "  -2.67%  " },
    @{ Row = 6; B = "Solana"; C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D = "'149.52"; E = "  -4.16%  " },
    @{ Row = 7; B = "LidoStakedEther"; C = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"; D = "3.504.63"; E = "  -2.46%  " },
    @{ Row = 8; B = "USDC"; C = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D = "'1.00"; E = "  -0.03%  " },
    @{ Row = 9; B = "XRP"; C = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D = "'0.481"; E = "  -1.49%  " },
    @{ Row = 10; B = "Dogecoin"; C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D = "'0.139"; E = "  -2.46%  " },
    @{ Row = 11; B = "Toncoin"; C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D = "'7.05"; E = "  +1.15%  " },
    @{ Row = 12; B = "Cardano"; C = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D = "'0.424"; E = "  -2.64%  " },
    @{ Row = 13; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "'0.0000219"; E = "  -3.04%  " },
    @{ Row = 14; B = "WrappedliquidstakedEther2.0"; C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D = "4.097.60"; E = "  -2.59%  " },
    @{ Row = 15; B = "Avalanche"; C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D = "'31.67"; E = "  -1.52%  " },
    @{ Row = 16; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "3.504.91"; E = "  -2.44%  " },
    @{ Row = 17; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "67.379.25"; E = "  -1.41%  " },
    @{ Row = 18; B = "TRON"; C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D = "'0.117"; E = "  -0.07%  " },
    @{ Row = 19; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "'6.41"; E = "  -0.59%  " },
    @{ Row = 20; B = "Chainlink"; C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "'15.16"; E = "  -3.56%  " },
    @{ Row = 21; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "'446.10"; E = "  -3.03%  " },
    @{ Row = 22; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "'9.22"; E = "  -6.23%  " },
    @{ Row = 23; B = "Polygon"; C = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D = "'0.623"; E = "  -3.35%  " },
    @{ Row = 24; B = "Litecoin"; C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "'77.45"; E = "  -0.98%  " },
    @{ Row = 25; B = "PEPE"; C = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D = "'0.0000128"; E = "  +9.42%  " },
    @{ Row = 26; B = "WrappedeETH"; C = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"; D = "3.644.47"; E = "  -2.56%  " },
    @{ Row = 27; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "'1.00"; E = "  +0.04%  " },
    @{ Row = 28; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "'10.27"; E = "  -4.59%  " },
    @{ Row = 29; B = "RenderToken"; C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D = "'8.33"; E = "  -1.16%  " },
    @{ Row = 30; B = "PancakeSwap"; C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D = "'2.49"; E = "  -3.91%  " },
    @{ Row = 31; B = "Binance-PegBSC-USD"; C = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"; D = "'1.00"; E = "  +0.05%  " },
    @{ Row = 32; B = "Fetch.AI"; C = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D = "'1.54"; E = "  -6.59%  " },
    @{ Row = 33; B = "Kaspa"; C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D = "'0.164"; E = "  +2.75%  " },
    @{ Row = 34; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "'25.66"; E = "  -1.96%  " },
    @{ Row = 35; B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "'6.10"; E = "  -1.72%  " },
    @{ Row = 36; B = "RenzoRestakedETH"; C = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"; D = "3.496.78"; E = "  -2.80%  " },
    @{ Row = 37; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "'1.83"; E = "  -4.77%  " },
    @{ Row = 38; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "'8.03"; E = "  -1.48%  " },
    @{ Row = 39; B = "USDe"; C = "https://coinranking.com/coin/exbfr2U-0+usde-usde"; D = "'1.00"; E = "  +0.01%  " },
    @{ Row = 40; B = "FirstDigitalUSD"; C = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D = "'1.00"; E = "  +0.06%  " },
    @{ Row = 41; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "'176.54"; E = "  -0.59%  " },
    @{ Row = 42; B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D = "'2.17"; E = "  +1.45%  " },
    @{ Row = 43; B = "Hedera"; C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "'0.0871"; E = "  -1.77%  " },
    @{ Row = 44; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "'5.41"; E = "  -4.22%  " },
    @{ Row = 45; B = "Mantle"; C = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D = "'0.878"; E = "  -2.73%  " },
    @{ Row = 46; B = "OKB"; C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D = "'45.38"; E = "  -1.57%  " },
    @{ Row = 47; B = "InjectiveProtocol"; C = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D = "'27.12"; E = "  -4.84%  " },
    @{ Row = 48; B = "ONDO"; C = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"; D = "'1.26"; E = "  +3.49%  " },
    @{ Row = 49; B = "dogwifhat"; C = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D = "'2.57"; E = "  -1.34%  " },
    @{ Row = 50; B = "Cosmos"; C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D = "'7.56"; E = "  -2.22%  " },
    @{ Row = 51; B = "SuiNetwork"; C = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"; D = "'0.993"; E = "  -2.15%  " }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
}

Write-Output "Updated $($rows.Count) rows"
